$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 39) that duplicates the last existing data row (row 38):
#   A39 -> same date text as A38 ("24-10-2025")
#   B39 -> same gold-price description text as B38
# New cells naturally pick up the column-level style (A=3, B=4), matching
# rows 33-38 which already follow this exact pattern.
$ws.Range("A39").Value2 = $ws.Range("A38").Value2
$ws.Range("B39").Value2 = $ws.Range("B38").Value2
